$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is a bare "NN%" string: force Text format first so Excel
# keeps the literal text instead of silently converting it to a numeric percentage.
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H34").NumberFormat = "@"

$ws.Range("E2").Value = "2026-02-07 08:17:57"
$ws.Range("E3").Value = "2026-02-07 08:18:00"
$ws.Range("N3").Value = "-7.7 °C 7:56 TU"
$ws.Range("O3").Value = "-6.1 °C"
$ws.Range("E4").Value = "2026-02-07 08:18:02"
$ws.Range("H4").Value = "59%"
$ws.Range("J4").Value = "1001.7 hPa"
$ws.Range("K4").Value = "0.1 MJ/m2"
$ws.Range("O4").Value = "10.9 °C"
$ws.Range("E5").Value = "2026-02-07 08:18:04"
$ws.Range("J5").Value = "1001.8 hPa"
$ws.Range("K5").Value = "0.1 MJ/m2"
$ws.Range("O5").Value = "8.1 °C"
$ws.Range("E6").Value = "2026-02-07 08:18:07"
$ws.Range("J6").Value = "1003.3 hPa"
$ws.Range("K6").Value = "0.1 MJ/m2"
$ws.Range("L6").Value = "38.9 km/h - 330º 7:50 TU"
$ws.Range("N6").Value = "11.2 °C 7:48 TU"
$ws.Range("O6").Value = "11.7 °C"
$ws.Range("E7").Value = "2026-02-07 08:18:10"
$ws.Range("H7").Value = "73%"
$ws.Range("J7").Value = "1003.0 hPa"
$ws.Range("K7").Value = "0.2 MJ/m2"
$ws.Range("O7").Value = "7.6 °C"
$ws.Range("E8").Value = "2026-02-07 08:18:13"
$ws.Range("K8").Value = "0.2 MJ/m2"
$ws.Range("M8").Value = "7.1 °C 7:57 TU"
$ws.Range("O8").Value = "4.0 °C"
$ws.Range("E9").Value = "2026-02-07 08:18:15"
$ws.Range("E10").Value = "2026-02-07 08:18:18"
$ws.Range("H10").Value = "99%"
$ws.Range("M10").Value = "9.7 °C 7:40 TU"
$ws.Range("O10").Value = "7.4 °C"
$ws.Range("E11").Value = "2026-02-07 08:18:21"
$ws.Range("E12").Value = "2026-02-07 08:18:23"
$ws.Range("O12").Value = "9.9 °C"
$ws.Range("E13").Value = "2026-02-07 08:18:25"
$ws.Range("H13").Value = "84%"
$ws.Range("M13").Value = "11.7 °C 7:57 TU"
$ws.Range("O13").Value = "8.2 °C"
$ws.Range("E14").Value = "2026-02-07 08:18:28"
$ws.Range("H14").Value = "71%"
$ws.Range("K14").Value = "0.0 MJ/m2"
$ws.Range("E15").Value = "2026-02-07 08:18:31"
$ws.Range("H15").Value = "90%"
$ws.Range("J15").Value = "1002.1 hPa"
$ws.Range("K15").Value = "0.1 MJ/m2"
$ws.Range("E16").Value = "2026-02-07 08:18:33"
$ws.Range("E17").Value = "2026-02-07 08:18:36"
$ws.Range("J17").Value = "1005.3 hPa"
$ws.Range("E18").Value = "2026-02-07 08:18:39"
$ws.Range("K18").Value = "0.1 MJ/m2"
$ws.Range("O18").Value = "-7.9 °C"
$ws.Range("E19").Value = "2026-02-07 08:18:41"
$ws.Range("J19").Value = "1006.7 hPa"
$ws.Range("K19").Value = "0.2 MJ/m2"
$ws.Range("L19").Value = "14.4 km/h - 271º 7:46 TU"
$ws.Range("E20").Value = "2026-02-07 08:18:44"
$ws.Range("H20").Value = "80%"
$ws.Range("K20").Value = "0.1 MJ/m2"
$ws.Range("E21").Value = "2026-02-07 08:18:46"
$ws.Range("J21").Value = "1002.5 hPa"
$ws.Range("K21").Value = "0.1 MJ/m2"
$ws.Range("O21").Value = "5.5 °C"
$ws.Range("E22").Value = "2026-02-07 08:18:49"
$ws.Range("K22").Value = "0.2 MJ/m2"
$ws.Range("O22").Value = "6.8 °C"
$ws.Range("E23").Value = "2026-02-07 08:18:51"
$ws.Range("J23").Value = "1001.9 hPa"
$ws.Range("K23").Value = "0.2 MJ/m2"
$ws.Range("L23").Value = "23.8 km/h - 31º 7:46 TU"
$ws.Range("M23").Value = "10.1 °C 7:59 TU"
$ws.Range("O23").Value = "7.5 °C"
$ws.Range("E24").Value = "2026-02-07 08:18:54"
$ws.Range("J24").Value = "1001.1 hPa"
$ws.Range("K24").Value = "0.2 MJ/m2"
$ws.Range("E25").Value = "2026-02-07 08:18:57"
$ws.Range("E26").Value = "2026-02-07 08:19:00"
$ws.Range("K26").Value = "0.1 MJ/m2"
$ws.Range("E27").Value = "2026-02-07 08:19:02"
$ws.Range("H27").Value = "90%"
$ws.Range("J27").Value = "1001.6 hPa"
$ws.Range("K27").Value = "0.1 MJ/m2"
$ws.Range("M27").Value = "12.0 °C 7:59 TU"
$ws.Range("O27").Value = "9.2 °C"
$ws.Range("E28").Value = "2026-02-07 08:19:04"
$ws.Range("J28").Value = "1004.5 hPa"
$ws.Range("N28").Value = "1.3 °C 7:33 TU"
$ws.Range("O28").Value = "2.8 °C"
$ws.Range("E29").Value = "2026-02-07 08:19:07"
$ws.Range("K29").Value = "0.1 MJ/m2"
$ws.Range("E30").Value = "2026-02-07 08:19:09"
$ws.Range("H30").Value = "78%"
$ws.Range("K30").Value = "0.2 MJ/m2"
$ws.Range("N30").Value = "-5.7 °C 7:59 TU"
$ws.Range("E31").Value = "2026-02-07 08:19:12"
$ws.Range("J31").Value = "1006.2 hPa"
$ws.Range("E32").Value = "2026-02-07 08:19:15"
$ws.Range("H32").Value = "56%"
$ws.Range("J32").Value = "1004.8 hPa"
$ws.Range("K32").Value = "0.1 MJ/m2"
$ws.Range("E33").Value = "2026-02-07 08:19:18"
$ws.Range("O33").Value = "6.8 °C"
$ws.Range("E34").Value = "2026-02-07 08:19:20"
$ws.Range("H34").Value = "80%"
$ws.Range("E35").Value = "2026-02-07 08:19:23"
$ws.Range("E36").Value = "2026-02-07 08:19:26"
$ws.Range("J36").Value = "1007.0 hPa"
$ws.Range("K36").Value = "0.1 MJ/m2"
$ws.Range("N36").Value = "2.4 °C 7:48 TU"
$ws.Range("O36").Value = "4.3 °C"
